# Refresh LR-pair data with new TPM-derived values (Sfrp1-Fzd2).
# Sending-cluster categories were recomputed/renamed:
#   old "Inflammatory-Mac" sending-cluster rows -> new "MuSCs" sending-cluster rows
#   old "MuSCs" sending-cluster rows            -> new "Neutrophils" sending-cluster rows
# All Ligand/Receptor/Edge metrics (columns E:T) are refreshed with the new TPM-based values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = 'ECs'
$ws.Cells.Item(2, 2).Value = 'Sfrp1'
$ws.Cells.Item(2, 3).Value = 'Fzd2'
$ws.Cells.Item(2, 4).Value = 'ECs'
$ws.Cells.Item(2, 5).Value = 2.0
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 0.2675626666666667
$ws.Cells.Item(2, 8).Value = 0.802688
$ws.Cells.Item(2, 9).Value = 0.01337792263105189
$ws.Cells.Item(2, 10).Value = 0.01337792263105189
$ws.Cells.Item(2, 11).Value = 3.0
$ws.Cells.Item(2, 12).Value = 1.0
$ws.Cells.Item(2, 13).Value = 0.3007906666666667
$ws.Cells.Item(2, 14).Value = 0.902372
$ws.Cells.Item(2, 15).Value = 0.03537029821880876
$ws.Cells.Item(2, 16).Value = 0.03537029821880876
$ws.Cells.Item(2, 17).Value = 0.08048035288177778
$ws.Cells.Item(2, 18).Value = 0.7243231759359999
$ws.Cells.Item(2, 19).Value = 0.0004731811130084563
$ws.Cells.Item(2, 20).Value = 0.0004731811130084562
# Row 3
$ws.Cells.Item(3, 1).Value = 'ECs'
$ws.Cells.Item(3, 2).Value = 'Sfrp1'
$ws.Cells.Item(3, 3).Value = 'Fzd2'
$ws.Cells.Item(3, 4).Value = 'FAPs'
$ws.Cells.Item(3, 5).Value = 2.0
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 0.2675626666666667
$ws.Cells.Item(3, 8).Value = 0.802688
$ws.Cells.Item(3, 9).Value = 0.01337792263105189
$ws.Cells.Item(3, 10).Value = 0.01337792263105189
$ws.Cells.Item(3, 11).Value = 3.0
$ws.Cells.Item(3, 12).Value = 1.0
$ws.Cells.Item(3, 13).Value = 7.647982
$ws.Cells.Item(3, 14).Value = 22.943946
$ws.Cells.Item(3, 15).Value = 0.899334434508434
$ws.Cells.Item(3, 16).Value = 0.899334434508434
$ws.Cells.Item(3, 17).Value = 2.046314458538667
$ws.Cells.Item(3, 18).Value = 18.416830126848
$ws.Cells.Item(3, 19).Value = 0.01203122648429464
$ws.Cells.Item(3, 20).Value = 0.01203122648429463
# Row 4
$ws.Cells.Item(4, 1).Value = 'ECs'
$ws.Cells.Item(4, 2).Value = 'Sfrp1'
$ws.Cells.Item(4, 3).Value = 'Fzd2'
$ws.Cells.Item(4, 4).Value = 'Inflammatory-Mac'
$ws.Cells.Item(4, 5).Value = 2.0
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 0.2675626666666667
$ws.Cells.Item(4, 8).Value = 0.802688
$ws.Cells.Item(4, 9).Value = 0.01337792263105189
$ws.Cells.Item(4, 10).Value = 0.01337792263105189
$ws.Cells.Item(4, 11).Value = 1.0
$ws.Cells.Item(4, 12).Value = 0.3333333333333333
$ws.Cells.Item(4, 13).Value = 0.0008990000000000001
$ws.Cells.Item(4, 14).Value = 0.002697
$ws.Cells.Item(4, 15).Value = 0.0001057143775473167
$ws.Cells.Item(4, 16).Value = 0.0001057143775473167
$ws.Cells.Item(4, 17).Value = 0.0002405388373333334
$ws.Cells.Item(4, 18).Value = 0.002164849536
$ws.Cells.Item(4, 19).Value = 0.000001414238763817812
$ws.Cells.Item(4, 20).Value = 0.000001414238763817812
# Row 5
$ws.Cells.Item(5, 1).Value = 'ECs'
$ws.Cells.Item(5, 2).Value = 'Sfrp1'
$ws.Cells.Item(5, 3).Value = 'Fzd2'
$ws.Cells.Item(5, 4).Value = 'MuSCs'
$ws.Cells.Item(5, 5).Value = 2.0
$ws.Cells.Item(5, 6).Value = 0.6666666666666666
$ws.Cells.Item(5, 7).Value = 0.2675626666666667
$ws.Cells.Item(5, 8).Value = 0.802688
$ws.Cells.Item(5, 9).Value = 0.01337792263105189
$ws.Cells.Item(5, 10).Value = 0.01337792263105189
$ws.Cells.Item(5, 11).Value = 3.0
$ws.Cells.Item(5, 12).Value = 1.0
$ws.Cells.Item(5, 13).Value = 0.5528646666666667
$ws.Cells.Item(5, 14).Value = 1.658594
$ws.Cells.Item(5, 15).Value = 0.06501195117304938
$ws.Cells.Item(5, 16).Value = 0.06501195117304936
$ws.Cells.Item(5, 17).Value = 0.1479259445191111
$ws.Cells.Item(5, 18).Value = 1.331333500672
$ws.Cells.Item(5, 19).Value = 0.000869724852886778
$ws.Cells.Item(5, 20).Value = 0.0008697248528867777
# Row 6
$ws.Cells.Item(6, 1).Value = 'ECs'
$ws.Cells.Item(6, 2).Value = 'Sfrp1'
$ws.Cells.Item(6, 3).Value = 'Fzd2'
$ws.Cells.Item(6, 4).Value = 'Resolving-Mac'
$ws.Cells.Item(6, 5).Value = 2.0
$ws.Cells.Item(6, 6).Value = 0.6666666666666666
$ws.Cells.Item(6, 7).Value = 0.2675626666666667
$ws.Cells.Item(6, 8).Value = 0.802688
$ws.Cells.Item(6, 9).Value = 0.01337792263105189
$ws.Cells.Item(6, 10).Value = 0.01337792263105189
$ws.Cells.Item(6, 11).Value = 1.0
$ws.Cells.Item(6, 12).Value = 0.3333333333333333
$ws.Cells.Item(6, 13).Value = 0.001510333333333333
$ws.Cells.Item(6, 14).Value = 0.004531
$ws.Cells.Item(6, 15).Value = 0.0001776017221605087
$ws.Cells.Item(6, 16).Value = 0.0001776017221605087
$ws.Cells.Item(6, 17).Value = 0.0004041088142222222
$ws.Cells.Item(6, 18).Value = 0.003636979328
$ws.Cells.Item(6, 19).Value = 0.00000237594209820486
$ws.Cells.Item(6, 20).Value = 0.00000237594209820486
# Row 7
$ws.Cells.Item(7, 1).Value = 'FAPs'
$ws.Cells.Item(7, 2).Value = 'Sfrp1'
$ws.Cells.Item(7, 3).Value = 'Fzd2'
$ws.Cells.Item(7, 4).Value = 'ECs'
$ws.Cells.Item(7, 5).Value = 3.0
$ws.Cells.Item(7, 6).Value = 1.0
$ws.Cells.Item(7, 7).Value = 18.95196233333333
$ws.Cells.Item(7, 8).Value = 56.855887
$ws.Cells.Item(7, 9).Value = 0.9475831922313891
$ws.Cells.Item(7, 10).Value = 0.947583192231389
$ws.Cells.Item(7, 11).Value = 3.0
$ws.Cells.Item(7, 12).Value = 1.0
$ws.Cells.Item(7, 13).Value = 0.3007906666666667
$ws.Cells.Item(7, 14).Value = 0.902372
$ws.Cells.Item(7, 15).Value = 0.03537029821880876
$ws.Cells.Item(7, 16).Value = 0.03537029821880876
$ws.Cells.Item(7, 17).Value = 5.700573384884889
$ws.Cells.Item(7, 18).Value = 51.305160463964
$ws.Cells.Item(7, 19).Value = 0.03351630009635503
$ws.Cells.Item(7, 20).Value = 0.03351630009635502
# Row 8
$ws.Cells.Item(8, 1).Value = 'FAPs'
$ws.Cells.Item(8, 2).Value = 'Sfrp1'
$ws.Cells.Item(8, 3).Value = 'Fzd2'
$ws.Cells.Item(8, 4).Value = 'FAPs'
$ws.Cells.Item(8, 5).Value = 3.0
$ws.Cells.Item(8, 6).Value = 1.0
$ws.Cells.Item(8, 7).Value = 18.95196233333333
$ws.Cells.Item(8, 8).Value = 56.855887
$ws.Cells.Item(8, 9).Value = 0.9475831922313891
$ws.Cells.Item(8, 10).Value = 0.947583192231389
$ws.Cells.Item(8, 11).Value = 3.0
$ws.Cells.Item(8, 12).Value = 1.0
$ws.Cells.Item(8, 13).Value = 7.647982
$ws.Cells.Item(8, 14).Value = 22.943946
$ws.Cells.Item(8, 15).Value = 0.899334434508434
$ws.Cells.Item(8, 16).Value = 0.899334434508434
$ws.Cells.Item(8, 17).Value = 144.9442667900113
$ws.Cells.Item(8, 18).Value = 1304.498401110102
$ws.Cells.Item(8, 19).Value = 0.852194194335113
$ws.Cells.Item(8, 20).Value = 0.8521941943351129
# Row 9
$ws.Cells.Item(9, 1).Value = 'FAPs'
$ws.Cells.Item(9, 2).Value = 'Sfrp1'
$ws.Cells.Item(9, 3).Value = 'Fzd2'
$ws.Cells.Item(9, 4).Value = 'Inflammatory-Mac'
$ws.Cells.Item(9, 5).Value = 3.0
$ws.Cells.Item(9, 6).Value = 1.0
$ws.Cells.Item(9, 7).Value = 18.95196233333333
$ws.Cells.Item(9, 8).Value = 56.855887
$ws.Cells.Item(9, 9).Value = 0.9475831922313891
$ws.Cells.Item(9, 10).Value = 0.947583192231389
$ws.Cells.Item(9, 11).Value = 1.0
$ws.Cells.Item(9, 12).Value = 0.3333333333333333
$ws.Cells.Item(9, 13).Value = 0.0008990000000000001
$ws.Cells.Item(9, 14).Value = 0.002697
$ws.Cells.Item(9, 15).Value = 0.0001057143775473167
$ws.Cells.Item(9, 16).Value = 0.0001057143775473167
$ws.Cells.Item(9, 17).Value = 0.01703781413766667
$ws.Cells.Item(9, 18).Value = 0.153340327239
$ws.Cells.Item(9, 19).Value = 0.0001001731673410406
$ws.Cells.Item(9, 20).Value = 0.0001001731673410406
# Row 10
$ws.Cells.Item(10, 1).Value = 'FAPs'
$ws.Cells.Item(10, 2).Value = 'Sfrp1'
$ws.Cells.Item(10, 3).Value = 'Fzd2'
$ws.Cells.Item(10, 4).Value = 'MuSCs'
$ws.Cells.Item(10, 5).Value = 3.0
$ws.Cells.Item(10, 6).Value = 1.0
$ws.Cells.Item(10, 7).Value = 18.95196233333333
$ws.Cells.Item(10, 8).Value = 56.855887
$ws.Cells.Item(10, 9).Value = 0.9475831922313891
$ws.Cells.Item(10, 10).Value = 0.947583192231389
$ws.Cells.Item(10, 11).Value = 3.0
$ws.Cells.Item(10, 12).Value = 1.0
$ws.Cells.Item(10, 13).Value = 0.5528646666666667
$ws.Cells.Item(10, 14).Value = 1.658594
$ws.Cells.Item(10, 15).Value = 0.06501195117304938
$ws.Cells.Item(10, 16).Value = 0.06501195117304936
$ws.Cells.Item(10, 17).Value = 10.47787033809756
$ws.Cells.Item(10, 18).Value = 94.300833042878
$ws.Cells.Item(10, 19).Value = 0.06160423222574933
$ws.Cells.Item(10, 20).Value = 0.06160423222574931
# Row 11
$ws.Cells.Item(11, 1).Value = 'FAPs'
$ws.Cells.Item(11, 2).Value = 'Sfrp1'
$ws.Cells.Item(11, 3).Value = 'Fzd2'
$ws.Cells.Item(11, 4).Value = 'Resolving-Mac'
$ws.Cells.Item(11, 5).Value = 3.0
$ws.Cells.Item(11, 6).Value = 1.0
$ws.Cells.Item(11, 7).Value = 18.95196233333333
$ws.Cells.Item(11, 8).Value = 56.855887
$ws.Cells.Item(11, 9).Value = 0.9475831922313891
$ws.Cells.Item(11, 10).Value = 0.947583192231389
$ws.Cells.Item(11, 11).Value = 1.0
$ws.Cells.Item(11, 12).Value = 0.3333333333333333
$ws.Cells.Item(11, 13).Value = 0.001510333333333333
$ws.Cells.Item(11, 14).Value = 0.004531
$ws.Cells.Item(11, 15).Value = 0.0001776017221605087
$ws.Cells.Item(11, 16).Value = 0.0001776017221605087
$ws.Cells.Item(11, 17).Value = 0.02862378044411111
$ws.Cells.Item(11, 18).Value = 0.257614023997
$ws.Cells.Item(11, 19).Value = 0.000168292406830647
$ws.Cells.Item(11, 20).Value = 0.000168292406830647
# Row 12
$ws.Cells.Item(12, 1).Value = 'MuSCs'
$ws.Cells.Item(12, 2).Value = 'Sfrp1'
$ws.Cells.Item(12, 3).Value = 'Fzd2'
$ws.Cells.Item(12, 4).Value = 'ECs'
$ws.Cells.Item(12, 5).Value = 3.0
$ws.Cells.Item(12, 6).Value = 1.0
$ws.Cells.Item(12, 7).Value = 0.7737046666666667
$ws.Cells.Item(12, 8).Value = 2.321114
$ws.Cells.Item(12, 9).Value = 0.03868462405050454
$ws.Cells.Item(12, 10).Value = 0.03868462405050453
$ws.Cells.Item(12, 11).Value = 3.0
$ws.Cells.Item(12, 12).Value = 1.0
$ws.Cells.Item(12, 13).Value = 0.3007906666666667
$ws.Cells.Item(12, 14).Value = 0.902372
$ws.Cells.Item(12, 15).Value = 0.03537029821880876
$ws.Cells.Item(12, 16).Value = 0.03537029821880876
$ws.Cells.Item(12, 17).Value = 0.2327231424897778
$ws.Cells.Item(12, 18).Value = 2.094508282408
$ws.Cells.Item(12, 19).Value = 0.001368286689148847
$ws.Cells.Item(12, 20).Value = 0.001368286689148847
# Row 13
$ws.Cells.Item(13, 1).Value = 'MuSCs'
$ws.Cells.Item(13, 2).Value = 'Sfrp1'
$ws.Cells.Item(13, 3).Value = 'Fzd2'
$ws.Cells.Item(13, 4).Value = 'FAPs'
$ws.Cells.Item(13, 5).Value = 3.0
$ws.Cells.Item(13, 6).Value = 1.0
$ws.Cells.Item(13, 7).Value = 0.7737046666666667
$ws.Cells.Item(13, 8).Value = 2.321114
$ws.Cells.Item(13, 9).Value = 0.03868462405050454
$ws.Cells.Item(13, 10).Value = 0.03868462405050453
$ws.Cells.Item(13, 11).Value = 3.0
$ws.Cells.Item(13, 12).Value = 1.0
$ws.Cells.Item(13, 13).Value = 7.647982
$ws.Cells.Item(13, 14).Value = 22.943946
$ws.Cells.Item(13, 15).Value = 0.899334434508434
$ws.Cells.Item(13, 16).Value = 0.899334434508434
$ws.Cells.Item(13, 17).Value = 5.917279363982667
$ws.Cells.Item(13, 18).Value = 53.25551427584401
$ws.Cells.Item(13, 19).Value = 0.03479041449463186
$ws.Cells.Item(13, 20).Value = 0.03479041449463185
# Row 14
$ws.Cells.Item(14, 1).Value = 'MuSCs'
$ws.Cells.Item(14, 2).Value = 'Sfrp1'
$ws.Cells.Item(14, 3).Value = 'Fzd2'
$ws.Cells.Item(14, 4).Value = 'Inflammatory-Mac'
$ws.Cells.Item(14, 5).Value = 3.0
$ws.Cells.Item(14, 6).Value = 1.0
$ws.Cells.Item(14, 7).Value = 0.7737046666666667
$ws.Cells.Item(14, 8).Value = 2.321114
$ws.Cells.Item(14, 9).Value = 0.03868462405050454
$ws.Cells.Item(14, 10).Value = 0.03868462405050453
$ws.Cells.Item(14, 11).Value = 1.0
$ws.Cells.Item(14, 12).Value = 0.3333333333333333
$ws.Cells.Item(14, 13).Value = 0.0008990000000000001
$ws.Cells.Item(14, 14).Value = 0.002697
$ws.Cells.Item(14, 15).Value = 0.0001057143775473167
$ws.Cells.Item(14, 16).Value = 0.0001057143775473167
$ws.Cells.Item(14, 17).Value = 0.0006955604953333334
$ws.Cells.Item(14, 18).Value = 0.006260044458000001
$ws.Cells.Item(14, 19).Value = 0.000004089520952151044
$ws.Cells.Item(14, 20).Value = 0.000004089520952151044
# Row 15
$ws.Cells.Item(15, 1).Value = 'MuSCs'
$ws.Cells.Item(15, 2).Value = 'Sfrp1'
$ws.Cells.Item(15, 3).Value = 'Fzd2'
$ws.Cells.Item(15, 4).Value = 'MuSCs'
$ws.Cells.Item(15, 5).Value = 3.0
$ws.Cells.Item(15, 6).Value = 1.0
$ws.Cells.Item(15, 7).Value = 0.7737046666666667
$ws.Cells.Item(15, 8).Value = 2.321114
$ws.Cells.Item(15, 9).Value = 0.03868462405050454
$ws.Cells.Item(15, 10).Value = 0.03868462405050453
$ws.Cells.Item(15, 11).Value = 3.0
$ws.Cells.Item(15, 12).Value = 1.0
$ws.Cells.Item(15, 13).Value = 0.5528646666666667
$ws.Cells.Item(15, 14).Value = 1.658594
$ws.Cells.Item(15, 15).Value = 0.06501195117304938
$ws.Cells.Item(15, 16).Value = 0.06501195117304936
$ws.Cells.Item(15, 17).Value = 0.4277539726351111
$ws.Cells.Item(15, 18).Value = 3.849785753716
$ws.Cells.Item(15, 19).Value = 0.002514962889919173
$ws.Cells.Item(15, 20).Value = 0.002514962889919172
# Row 16
$ws.Cells.Item(16, 1).Value = 'MuSCs'
$ws.Cells.Item(16, 2).Value = 'Sfrp1'
$ws.Cells.Item(16, 3).Value = 'Fzd2'
$ws.Cells.Item(16, 4).Value = 'Resolving-Mac'
$ws.Cells.Item(16, 5).Value = 3.0
$ws.Cells.Item(16, 6).Value = 1.0
$ws.Cells.Item(16, 7).Value = 0.7737046666666667
$ws.Cells.Item(16, 8).Value = 2.321114
$ws.Cells.Item(16, 9).Value = 0.03868462405050454
$ws.Cells.Item(16, 10).Value = 0.03868462405050453
$ws.Cells.Item(16, 11).Value = 1.0
$ws.Cells.Item(16, 12).Value = 0.3333333333333333
$ws.Cells.Item(16, 13).Value = 0.001510333333333333
$ws.Cells.Item(16, 14).Value = 0.004531
$ws.Cells.Item(16, 15).Value = 0.0001776017221605087
$ws.Cells.Item(16, 16).Value = 0.0001776017221605087
$ws.Cells.Item(16, 17).Value = 0.001168551948222222
$ws.Cells.Item(16, 18).Value = 0.010516967534
$ws.Cells.Item(16, 19).Value = 0.000006870455852501438
$ws.Cells.Item(16, 20).Value = 0.000006870455852501438
# Row 17
$ws.Cells.Item(17, 1).Value = 'Neutrophils'
$ws.Cells.Item(17, 2).Value = 'Sfrp1'
$ws.Cells.Item(17, 3).Value = 'Fzd2'
$ws.Cells.Item(17, 4).Value = 'ECs'
$ws.Cells.Item(17, 5).Value = 1.0
$ws.Cells.Item(17, 6).Value = 0.3333333333333333
$ws.Cells.Item(17, 7).Value = 0.007085333333333333
$ws.Cells.Item(17, 8).Value = 0.021256
$ws.Cells.Item(17, 9).Value = 0.0003542610870545456
$ws.Cells.Item(17, 10).Value = 0.0003542610870545455
$ws.Cells.Item(17, 11).Value = 3.0
$ws.Cells.Item(17, 12).Value = 1.0
$ws.Cells.Item(17, 13).Value = 0.3007906666666667
$ws.Cells.Item(17, 14).Value = 0.902372
$ws.Cells.Item(17, 15).Value = 0.03537029821880876
$ws.Cells.Item(17, 16).Value = 0.03537029821880876
$ws.Cells.Item(17, 17).Value = 0.002131202136888889
$ws.Cells.Item(17, 18).Value = 0.019180819232
$ws.Cells.Item(17, 19).Value = 0.00001253032029643865
$ws.Cells.Item(17, 20).Value = 0.00001253032029643865
# Row 18
$ws.Cells.Item(18, 1).Value = 'Neutrophils'
$ws.Cells.Item(18, 2).Value = 'Sfrp1'
$ws.Cells.Item(18, 3).Value = 'Fzd2'
$ws.Cells.Item(18, 4).Value = 'FAPs'
$ws.Cells.Item(18, 5).Value = 1.0
$ws.Cells.Item(18, 6).Value = 0.3333333333333333
$ws.Cells.Item(18, 7).Value = 0.007085333333333333
$ws.Cells.Item(18, 8).Value = 0.021256
$ws.Cells.Item(18, 9).Value = 0.0003542610870545456
$ws.Cells.Item(18, 10).Value = 0.0003542610870545455
$ws.Cells.Item(18, 11).Value = 3.0
$ws.Cells.Item(18, 12).Value = 1.0
$ws.Cells.Item(18, 13).Value = 7.647982
$ws.Cells.Item(18, 14).Value = 22.943946
$ws.Cells.Item(18, 15).Value = 0.899334434508434
$ws.Cells.Item(18, 16).Value = 0.899334434508434
$ws.Cells.Item(18, 17).Value = 0.05418850179733333
$ws.Cells.Item(18, 18).Value = 0.487696516176
$ws.Cells.Item(18, 19).Value = 0.0003185991943945429
$ws.Cells.Item(18, 20).Value = 0.0003185991943945428
# Row 19
$ws.Cells.Item(19, 1).Value = 'Neutrophils'
$ws.Cells.Item(19, 2).Value = 'Sfrp1'
$ws.Cells.Item(19, 3).Value = 'Fzd2'
$ws.Cells.Item(19, 4).Value = 'Inflammatory-Mac'
$ws.Cells.Item(19, 5).Value = 1.0
$ws.Cells.Item(19, 6).Value = 0.3333333333333333
$ws.Cells.Item(19, 7).Value = 0.007085333333333333
$ws.Cells.Item(19, 8).Value = 0.021256
$ws.Cells.Item(19, 9).Value = 0.0003542610870545456
$ws.Cells.Item(19, 10).Value = 0.0003542610870545455
$ws.Cells.Item(19, 11).Value = 1.0
$ws.Cells.Item(19, 12).Value = 0.3333333333333333
$ws.Cells.Item(19, 13).Value = 0.0008990000000000001
$ws.Cells.Item(19, 14).Value = 0.002697
$ws.Cells.Item(19, 15).Value = 0.0001057143775473167
$ws.Cells.Item(19, 16).Value = 0.0001057143775473167
$ws.Cells.Item(19, 17).Value = 0.000006369714666666667
$ws.Cells.Item(19, 18).Value = 0.000057327432
$ws.Cells.Item(19, 19).Value = 0.00000003745049030720706
$ws.Cells.Item(19, 20).Value = 0.00000003745049030720705
# Row 20
$ws.Cells.Item(20, 1).Value = 'Neutrophils'
$ws.Cells.Item(20, 2).Value = 'Sfrp1'
$ws.Cells.Item(20, 3).Value = 'Fzd2'
$ws.Cells.Item(20, 4).Value = 'MuSCs'
$ws.Cells.Item(20, 5).Value = 1.0
$ws.Cells.Item(20, 6).Value = 0.3333333333333333
$ws.Cells.Item(20, 7).Value = 0.007085333333333333
$ws.Cells.Item(20, 8).Value = 0.021256
$ws.Cells.Item(20, 9).Value = 0.0003542610870545456
$ws.Cells.Item(20, 10).Value = 0.0003542610870545455
$ws.Cells.Item(20, 11).Value = 3.0
$ws.Cells.Item(20, 12).Value = 1.0
$ws.Cells.Item(20, 13).Value = 0.5528646666666667
$ws.Cells.Item(20, 14).Value = 1.658594
$ws.Cells.Item(20, 15).Value = 0.06501195117304938
$ws.Cells.Item(20, 16).Value = 0.06501195117304936
$ws.Cells.Item(20, 17).Value = 0.003917230451555556
$ws.Cells.Item(20, 18).Value = 0.035255074064
$ws.Cells.Item(20, 19).Value = 0.00002303120449410151
$ws.Cells.Item(20, 20).Value = 0.0000230312044941015
# Row 21
$ws.Cells.Item(21, 1).Value = 'Neutrophils'
$ws.Cells.Item(21, 2).Value = 'Sfrp1'
$ws.Cells.Item(21, 3).Value = 'Fzd2'
$ws.Cells.Item(21, 4).Value = 'Resolving-Mac'
$ws.Cells.Item(21, 5).Value = 1.0
$ws.Cells.Item(21, 6).Value = 0.3333333333333333
$ws.Cells.Item(21, 7).Value = 0.007085333333333333
$ws.Cells.Item(21, 8).Value = 0.021256
$ws.Cells.Item(21, 9).Value = 0.0003542610870545456
$ws.Cells.Item(21, 10).Value = 0.0003542610870545455
$ws.Cells.Item(21, 11).Value = 1.0
$ws.Cells.Item(21, 12).Value = 0.3333333333333333
$ws.Cells.Item(21, 13).Value = 0.001510333333333333
$ws.Cells.Item(21, 14).Value = 0.004531
$ws.Cells.Item(21, 15).Value = 0.0001776017221605087
$ws.Cells.Item(21, 16).Value = 0.0001776017221605087
$ws.Cells.Item(21, 17).Value = 0.00001070121511111111
$ws.Cells.Item(21, 18).Value = 0.000096310936
$ws.Cells.Item(21, 19).Value = 0.00000006291737915534118
$ws.Cells.Item(21, 20).Value = 0.00000006291737915534118

Write-Host "Updated Sfrp1-Fzd2 LR-pair sheet with refreshed TPM values."
